$wb = $excel.ActiveWorkbook

# --- Summary sheet: update Total Trades and Win Rate % ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 49
$summary.Range("B9").Value = 42.86

# --- Strategy Status sheet: update MarketMaking row (Trades, Win Rate %) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 49
$status.Range("G4").Value = 42.86

# --- Append new closed trade (Trade #49) to both "All Trades" and
#     "MarketMaking" sheets, which mirror the same trade log. ---
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 50

    # Column B holds a "YYYY-MM-DD" string. Assigning such a string directly
    # causes the engine to auto-convert it into a date serial value, so
    # instead copy it from an existing cell that already holds the exact
    # same literal date text, preserving it as plain text.
    $ws.Cells.Item(2, 2).Copy($ws.Cells.Item($row, 2))

    $ws.Cells.Item($row, 1).Value = 49
    $ws.Cells.Item($row, 3).Value = "12:48:49"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.04
    $ws.Cells.Item($row, 7).Value = 0.04
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.16
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
